$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$nl = [char]10
$newText = "Conversión del día 💰" + $nl + `
    "✅ Dólar paralelo: 68" + $nl + `
    $nl + `
    "Binance" + $nl + `
    "✅ 1000 Bs = 7.14 = 28611.99 pesos" + $nl + `
    "✅ 28611.99 pesos = 7.11 = 969.8 Bs" + $nl + `
    $nl + `
    "Promedio competencia" + $nl + `
    "✅ Tasa pesos: 20" + $nl + `
    "✅ Tasa Bs: 20" + $nl + `
    "✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate cells N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 139.999
$ws2.Range("O10").Value = 4005.65
$ws2.Range("N12").Value = 4026.99
$ws2.Range("O12").Value = 136.495
